$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '43.286.80'
$ws.Range('E2').Value = '  -1.96%  '
$ws.Range('D3').Value = '2.237.71'
$ws.Range('E3').Value = '  -1.74%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = "'230.49"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.80%  '
$ws.Range('D6').Value = "'0.638"
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Value = "'63.47"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.59%  '
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('E9').Value = '  +1.09%  '
$ws.Range('D10').Value = "'0.0947"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -9.26%  '
$ws.Range('E11').Value = '  -0.13%  '
$ws.Range('D12').Value = "'27.59"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.63%  '
$ws.Range('E13').Value = '  -1.60%  '
$ws.Range('D14').Value = '2.570.28'
$ws.Range('E14').Value = '  -1.96%  '
$ws.Range('D15').Value = "'15.26"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.83%  '
$ws.Range('D16').Value = "'6.03"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.54%  '
$ws.Range('D17').Value = "'0.823"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('D18').Value = '2.243.31'
$ws.Range('E18').Value = '  -1.78%  '
$ws.Range('D19').Value = '43.173.45'
$ws.Range('E19').Value = '  -1.92%  '
$ws.Range('D20').Value = '0.0₃0961'
$ws.Range('E20').Value = '  -9.97%  '
$ws.Range('D21').Value = "'72.81"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.47%  '
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('D23').Value = "'245.64"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.15%  '
$ws.Range('D24').Value = "'1.00"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').Value = "'3.66"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +28.47%  '
$ws.Range('D26').Value = "'2.40"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.17%  '
$ws.Range('D27').Value = "'2.28"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.57%  '
$ws.Range('D28').Value = "'9.70"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.80%  '
$ws.Range('D29').Value = "'173.28"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.89%  '
$ws.Range('D30').Value = "'21.48"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.30%  '
$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').Value = "'0.128"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.59%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = "'1.40"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.28%  '
$ws.Range('E33').Value = '  -0.19%  '
$ws.Range('D34').Value = "'4.92"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.64%  '
$ws.Range('E35').Value = '  -1.62%  '
$ws.Range('D36').Value = "'4.89"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.16%  '
$ws.Range('E37').Value = '  -7.08%  '
$ws.Range('D38').Value = "'6.26"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.86%  '
$ws.Range('D39').Value = "'2.26"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.67%  '
$ws.Range('E40').Value = '  +0.15%  '
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('D42').Value = "'8.62"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.39%  '
$ws.Range('D43').Value = "'4.45"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.11%  '
$ws.Range('E44').Value = '  -4.23%  '
$ws.Range('D45').Value = "'96.24"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.88%  '
$ws.Range('D46').Value = "'0.0938"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.33%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').Value = "'1.18"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.30%  '
$ws.Range('B48').Value = 'TerraClassic'
$ws.Range('C48').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D48').Value = "'0.000207"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.48%  '
$ws.Range('D49').Value = '1.441.47'
$ws.Range('E49').Value = '  -2.00%  '
$ws.Range('D50').Value = "'9.81"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.94%  '
$ws.Range('B51').Value = 'HuobiToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D51').Value = "'2.74"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.66%  '
